$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = 'nicest-2-variables:10000'
$ws.Range("B20").Value = 'discharge'
$ws.Range("C20").Value = 'liquid water which drains from land'
$ws.Range("D20").Value = 'streamflow, runoff'

$ws.Range("A21").Value = 'nicest-2-variables:10001'
$ws.Range("B21").Value = 'aerosol optical thickness'
$ws.Range("C21").Value = ""

$ws.Range("A22").Value = 'nicest-2-variables:10002'
$ws.Range("B22").Value = 'fine mode optical thickness'
$ws.Range("C22").Value = ""

$ws.Range("A23").Value = 'nicest-2-variables:10003'
$ws.Range("B23").Value = 'coarse mode optical thickness'
$ws.Range("C23").Value = ""

$ws.Range("A24").Value = 'nicest-2-variables:10004'
$ws.Range("B24").Value = 'precipitation'
$ws.Range("C24").Value = ""

$ws.Range("A25").Value = 'nicest-2-variables:10005'
$ws.Range("B25").Value = 'relative humidity'
$ws.Range("C25").Value = 'indicates a present state of absolute humidity relative to a maximum humidity given the same temperature'

$ws.Range("A26").Value = 'nicest-2-variables:10006'
$ws.Range("B26").Value = 'global radiation'
$ws.Range("C26").Value = ""

$ws.Range("A27").Value = 'nicest-2-variables:10007'
$ws.Range("B27").Value = 'upper-air temperature'
$ws.Range("C27").Value = ""

$ws.Range("A28").Value = 'nicest-2-variables:10008'
$ws.Range("B28").Value = 'zonal/eastward wind'
$ws.Range("C28").Value = 'zonal wind (positive in a eastward direction)'

$ws.Range("A29").Value = 'nicest-2-variables:10009'
$ws.Range("B29").Value = 'meridional/northward wind'
$ws.Range("C29").Value = 'northward component of the near surface wind'

$ws.Range("A30").Value = 'nicest-2-variables:10010'
$ws.Range("B30").Value = 'vertical velocity/wind'
$ws.Range("C30").Value = ""

$ws.Range("A31").Value = 'nicest-2-variables:10011'
$ws.Range("B31").Value = 'sea ice thickness'
$ws.Range("C31").Value = 'actual (floe) thickness of sea ice'

$ws.Range("A32").Value = 'nicest-2-variables:10012'
$ws.Range("B32").Value = 'cloud fraction'
$ws.Range("C32").Value = ""

$ws.Range("A33").Value = 'nicest-2-variables:10013'
$ws.Range("B33").Value = 'ice fraction'
$ws.Range("C33").Value = ""

$ws.Range("A34").Value = 'nicest-2-variables:10014'
$ws.Range("B34").Value = 'atmosphere mass content of cloud liquid water'
$ws.Range("C34").Value = ""

$ws.Range("A35").Value = 'nicest-2-variables:10015'
$ws.Range("B35").Value = 'atmosphere mass content of cloud ice'
$ws.Range("C35").Value = ""

$ws.Range("A36").Value = 'nicest-2-variables:10016'
$ws.Range("B36").Value = 'ocean mixed layer depth'
$ws.Range("C36").Value = ""

$ws.Range("A37").Value = 'nicest-2-variables:10017'
$ws.Range("B37").Value = 'number of processor-hour per simulated year'
$ws.Range("C37").Value = ""

$ws.Range("A38").Value = 'nicest-2-variables:10018'
$ws.Range("B38").Value = 'number of simulated year per wall-clock day'
$ws.Range("C38").Value = ""

$ws.Range("A39").Value = 'nicest-2-variables:10019'
$ws.Range("B39").Value = 'number of computation seconds per model-day'
$ws.Range("C39").Value = ""

$ws.Range("A40").Value = 'nicest-2-variables:10020'
$ws.Range("B40").Value = 'dissolved inorganic carbon'
$ws.Range("C40").Value = ""

$ws.Rows("41:42").Delete()
